# Fruta / hortaliza, semanal
# The data rows (2-27, columns A-T) get rearranged into a new row order.
# Build a mapping from the NEW row number -> the OLD row number whose
# data should end up there, then snapshot all old values first (since
# this is a full permutation, not a simple shift) and write them back.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 2
$lastRow = 27
$lastCol = 20   # column T

# new row -> source (old) row
$rowMap = @{
    2  = 4
    3  = 5
    4  = 27
    5  = 20
    6  = 6
    7  = 7
    8  = 3
    9  = 21
    10 = 9
    11 = 8
    12 = 13
    13 = 25
    14 = 10
    15 = 12
    16 = 19
    17 = 16
    18 = 18
    19 = 14
    20 = 17
    21 = 23
    22 = 24
    23 = 26
    24 = 22
    25 = 15
    26 = 2
    27 = 11
}

# Snapshot every cell's current value before writing anything, since
# several source rows are reused as destinations for other rows.
$snapshot = @{}
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $rowVals = @()
    for ($c = 1; $c -le $lastCol; $c++) {
        $rowVals += , ($ws.Cells.Item($r, $c).Value2)
    }
    $snapshot[$r] = $rowVals
}

# Write back according to the mapping.
for ($newRow = $firstRow; $newRow -le $lastRow; $newRow++) {
    $srcRow = $rowMap[$newRow]
    $vals = $snapshot[$srcRow]
    for ($c = 1; $c -le $lastCol; $c++) {
        $ws.Cells.Item($newRow, $c).Value2 = $vals[$c - 1]
    }
}
